# Applies the "added new milestones and my role to ppt" edit.
$p = $ppt.ActivePresentation

# --- Slide 3: "Team Members & Responsibilities" ---
# Brendan Madigan's bullet gains his role description.
$s3 = $p.Slides.Item(3)
$shBrendan = $s3.Shapes.Item(8)
$shBrendan.TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Brendan Madigan – Developer: tool research and testing"

# --- Slide 4: "Technology Stack & Toolchain" ---
# Replace the last bullet (about problem solving) with a new bullet describing
# the GUI library (Tkinter) the team is using.
$s4 = $p.Slides.Item(4)
$shTech = $s4.Shapes.Item(5)
$tr4 = $shTech.TextFrame.TextRange
$lastPara = $tr4.Paragraphs(3)
$lastPara.Text = "For the GUI we are using the "
$runTkinter = $lastPara.InsertAfter("Tkinter")
$runTkinter.InsertAfter(" library") | Out-Null

# --- Slide 5: "Potential Future Milestones" ---
$s5 = $p.Slides.Item(5)
$shMilestones = $s5.Shapes.Item(7)
$trM = $shMilestones.TextFrame.TextRange
$trM.Paragraphs(1).Runs(1).Text = "Milestone 3 hopes to be able to change the name and directory of WAV files."
$trM.Paragraphs(2).Runs(1).Text = "Milestone 4 hopes to be able to visualize the WAV file."

# --- Slide 7: "Communication & Workflow Plan" ---
$s7 = $p.Slides.Item(7)
$shComm = $s7.Shapes.Item(5)
$trC = $shComm.TextFrame.TextRange
$trC.Paragraphs(3).Runs(1).Text = "Hunter Hutchison is the Maintainer, while the other four are developers"
